# "SE BORRO LA ULTIMA TABLA" - the last image ("Imagen 4") anchored in the
# final paragraph of the document is removed; Word leaves behind the
# paragraph mark together with the implicit "_GoBack" bookmark that Word
# drops at the position of the most recent edit.

$d = $word.ActiveDocument

# Locate and remove the last picture ("Imagen 4") anchored in the document.
for ($i = $d.Shapes.Count; $i -ge 1; $i--) {
    $shp = $d.Shapes.Item($i)
    if ($shp.Name -eq "Imagen 4") {
        $shp.Delete()
    }
}

# Word stamps a "_GoBack" bookmark at the site of the last edit (here, the
# now-empty trailing paragraph left behind once the picture was removed).
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$editRange = $lastPara.Range
$editRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $editRange)
